$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.5028904953788366
$ws.Range("D2").Value = 0.4865946797507227

$ws.Range("C3").Value = 0.1774724807286316
$ws.Range("D3").Value = 0.2420119890740031

$ws.Range("C4").Value = 0.06883397053823591
$ws.Range("D4").Value = 0.06585587115020564

$ws.Range("C5").Value = 0.03308301932852964
$ws.Range("D5").Value = 0.02839299407078133

$ws.Range("C6").Value = 0.01915927998178294
$ws.Range("D6").Value = 0.01625385792084194
